$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4942.0713
$ws.Range("I15").Value = 4942.0713
$ws.Range("K15").Value = 14826.2139
$ws.Range("M15").Value = -14657.2139

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 12230.318
$ws.Range("I137").Value = 19472.285
$ws.Range("J137").Value = 8850.733
$ws.Range("K137").Value = 58416.855
$ws.Range("L137").Value = 26552.199
$ws.Range("M137").Value = -55866.855
$ws.Range("N137").Value = -31652.199

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2858.1467
$ws.Range("I138").Value = 1811.9286
$ws.Range("J138").Value = 3481.4255
$ws.Range("K138").Value = 5435.7858
$ws.Range("L138").Value = 10444.2765
$ws.Range("M138").Value = -295.7857999999997
$ws.Range("N138").Value = -20724.2765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 37800
$ws.Range("J44").Value = 37800
$ws.Range("L44").Value = 37800
$ws.Range("N44").Value = -38776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3046.5
$ws.Range("I61").Value = 2699.5
$ws.Range("J61").Value = 3145.6428
$ws.Range("K61").Value = 2699.5
$ws.Range("L61").Value = 3145.6428
$ws.Range("M61").Value = -2487.5
$ws.Range("N61").Value = -3569.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1333.48
$ws.Range("I74").Value = 1190.3064
$ws.Range("K74").Value = 1190.3064
$ws.Range("M74").Value = -316.3063999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1333.48
$ws.Range("I77").Value = 1190.3064
$ws.Range("K77").Value = 5951.531999999999
$ws.Range("M77").Value = -1583.531999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 17860718
$ws.Range("I132").Value = 22730640
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 68191920
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -68189390
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3046.5
$ws.Range("I136").Value = 2699.5
$ws.Range("J136").Value = 3145.6428
$ws.Range("K136").Value = 8098.5
$ws.Range("L136").Value = 9436.928400000001
$ws.Range("M136").Value = -5548.5
$ws.Range("N136").Value = -14536.9284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3204.2144
$ws.Range("I134").Value = 2930.0715
$ws.Range("J134").Value = 3478.3572
$ws.Range("K134").Value = 8790.2145
$ws.Range("L134").Value = 10435.0716
$ws.Range("M134").Value = -6255.2145
$ws.Range("N134").Value = -15505.0716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2162.9333
$ws.Range("I58").Value = 1494.9474
$ws.Range("J58").Value = 3316.7273
$ws.Range("K58").Value = 1494.9474
$ws.Range("L58").Value = 3316.7273
$ws.Range("M58").Value = -1291.9474
$ws.Range("N58").Value = -3722.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 521630.03
$ws.Range("I132").Value = 2042.2106
$ws.Range("J132").Value = 1755651.1
$ws.Range("K132").Value = 6126.6318
$ws.Range("L132").Value = 5266953.300000001
$ws.Range("M132").Value = -3596.6318
$ws.Range("N132").Value = -5272013.300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 668665.9
$ws.Range("I134").Value = 1548.5
$ws.Range("J134").Value = 1275136.1
$ws.Range("K134").Value = 4645.5
$ws.Range("L134").Value = 3825408.3
$ws.Range("M134").Value = -2110.5
$ws.Range("N134").Value = -3830478.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2162.9333
$ws.Range("I136").Value = 1494.9474
$ws.Range("J136").Value = 3316.7273
$ws.Range("K136").Value = 4484.8422
$ws.Range("L136").Value = 9950.1819
$ws.Range("M136").Value = -1934.8422
$ws.Range("N136").Value = -15050.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2259.7742
$ws.Range("I68").Value = 1791.5834
$ws.Range("J68").Value = 2372.14
$ws.Range("K68").Value = 5374.7502
$ws.Range("L68").Value = 7116.42
$ws.Range("M68").Value = -4563.7502
$ws.Range("N68").Value = -8738.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2259.7742
$ws.Range("I71").Value = 1791.5834
$ws.Range("J71").Value = 2372.14
$ws.Range("K71").Value = 16124.2506
$ws.Range("L71").Value = 21349.26
$ws.Range("M71").Value = -12068.2506
$ws.Range("N71").Value = -29461.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 27780026
$ws.Range("I87").Value = 2997
$ws.Range("J87").Value = 111111110
$ws.Range("K87").Value = 8991
$ws.Range("L87").Value = 333333330
$ws.Range("M87").Value = -7743
$ws.Range("N87").Value = -333335826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 27780026
$ws.Range("I90").Value = 2997
$ws.Range("J90").Value = 111111110
$ws.Range("K90").Value = 26973
$ws.Range("L90").Value = 999999990
$ws.Range("M90").Value = -20733
$ws.Range("N90").Value = -1000012470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 14772
$ws.Range("I107").Value = 10390.4
$ws.Range("K107").Value = 31171.2
$ws.Range("M107").Value = -29251.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 43485110
$ws.Range("I137").Value = 2870
$ws.Range("J137").Value = 111124150
$ws.Range("K137").Value = 8610
$ws.Range("L137").Value = 333372450
$ws.Range("M137").Value = -3510
$ws.Range("N137").Value = -333382650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3093.7837
$ws.Range("I132").Value = 2641.0417
$ws.Range("J132").Value = 3929.6155
$ws.Range("K132").Value = 7923.125100000001
$ws.Range("L132").Value = 11788.8465
$ws.Range("M132").Value = -5393.125100000001
$ws.Range("N132").Value = -16848.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2346.2173
$ws.Range("I136").Value = 1632.3889
$ws.Range("J136").Value = 4916
$ws.Range("K136").Value = 4897.1667
$ws.Range("L136").Value = 14748
$ws.Range("M136").Value = -2347.1667
$ws.Range("N136").Value = -19848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1836.2979
$ws.Range("I132").Value = 1388.2433
$ws.Range("J132").Value = 3494.1
$ws.Range("K132").Value = 4164.7299
$ws.Range("L132").Value = 10482.3
$ws.Range("M132").Value = -1634.7299
$ws.Range("N132").Value = -15542.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 418108.34
$ws.Range("I136").Value = 715471.6
$ws.Range("J136").Value = 1799.7
$ws.Range("K136").Value = 2146414.8
$ws.Range("L136").Value = 5399.1
$ws.Range("M136").Value = -2143864.8
$ws.Range("N136").Value = -10499.1
